# Append new broker-coverage rows for DataDate 20210119 (hait_ehfz indirect method update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking date strings ("20210119"); temporarily format
# those cells as Text so typing the value doesn't get auto-converted to a
# number, matching the existing "DataDate" column's text storage.
$ws.Range("A27:A31").NumberFormat = "@"

$ws.Cells.Item(27, 1).Value = "20210119"
$ws.Cells.Item(27, 2).Value = "hait"
$ws.Cells.Item(27, 3).Value = "DataFileNotExists"

$ws.Cells.Item(28, 1).Value = "20210119"
$ws.Cells.Item(28, 2).Value = "huat"
$ws.Cells.Item(28, 3).Value = "DataFileNotExists"

$ws.Cells.Item(29, 1).Value = "20210119"
$ws.Cells.Item(29, 2).Value = "swhy"
$ws.Cells.Item(29, 3).Value = "DataFileNotExists"

$ws.Cells.Item(30, 1).Value = "20210119"
$ws.Cells.Item(30, 2).Value = "gtja"
$ws.Cells.Item(30, 3).Value = "DataFileNotExists"

$ws.Cells.Item(31, 1).Value = "20210119"
$ws.Cells.Item(31, 2).Value = "zx"
$ws.Cells.Item(31, 3).Value = 37

# Drop the temporary text number-format so the new cells keep the same
# (default) style as the rest of the table.
$ws.Range("A27:A31").ClearFormats()
